$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.69%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.65%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.867"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "14.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08131"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.03%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.616"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.21%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.779"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.84%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.979"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9477"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1327"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "13.57%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.2000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.22%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.948"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "42.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09362"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.28%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03502"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "6.75%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09642"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.80%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001315"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.62%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006376"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.44%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3542"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.82%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1408"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.03%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2420"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.10%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04441"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.94%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001264"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "5.87%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004454"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.14%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001091"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.32%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004011"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02454"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "13.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05285"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.76%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007558"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.03%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1436"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.69%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009067"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.45%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002053"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.54%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01052"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "35.31%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006838"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.87%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.83%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003513"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "23.47%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001810"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "7.46%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002111"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.83%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002010"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.83%"
